$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# MyForecast (column D) flag updates
$ws1.Range("D2").Value = 1
$ws1.Range("D3").Value = 1
$ws1.Range("D12").Value = 1
$ws1.Range("D13").Value = 1

# Seasonality Index (column L) updates
$ws1.Range("L3").Value = 0.95
$ws1.Range("L4").Value = 0.8100000000000001
$ws1.Range("L5").Value = 1.05
$ws1.Range("L6").Value = 0.97
$ws1.Range("L7").Value = 1.06
$ws1.Range("L8").Value = 1.07
$ws1.Range("L9").Value = 0.8100000000000001
$ws1.Range("L10").Value = 0.82
$ws1.Range("L11").Value = 1.16
$ws1.Range("L12").Value = 0.84
$ws1.Range("L13").Value = 0.83
$ws1.Range("L14").Value = 0.97
$ws1.Range("L15").Value = 1.06
$ws1.Range("L16").Value = 1.17
$ws1.Range("L17").Value = 0.93

# --- Sheet "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# Use a leading apostrophe so Excel keeps these as text values (matching
# the original cell type) instead of auto-converting to numbers.
$ws2.Range("B9").Value = "'10"
$ws2.Range("B10").Value = "'6"
$ws2.Range("B11").Value = "'4"
$ws2.Range("B14").Value = "'0"
